# Add a new "2020" column (L) that mirrors the existing "2019" column (K),
# then move the active selection to N5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: year 2020 (same style as the 2019 header cell, K4).
$ws.Range("L4").Value = 2020

# Copy the K5:K12 data values into L5:L12 (same figures as 2019).
$ws.Range("L5").Value = 5.6
$ws.Range("L6").Value = 0.8
$ws.Range("L7").Value = 1.9
$ws.Range("L8").Value = 0.7
$ws.Range("L9").Value = 0.7
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L12").Value = 0.2

# Move the selection, like the author did in the recorded session.
$ws.Range("N5").Select()
